$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.842.19"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "2.447.70"
$ws.Range("E3").Value = "  -3.03%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.72"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.69"
$ws.Range("E6").Value = "  -1.97%  "
$ws.Range("E7").Value = "  +0.49%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.562"
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("D9").Value = "2.453.31"
$ws.Range("E9").Value = "  -2.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0982"
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("E11").Value = "  -2.07%  "
$ws.Range("E12").Value = "  -4.72%  "
$ws.Range("E13").Value = "  -2.78%  "
$ws.Range("D14").Value = "2.883.17"
$ws.Range("E14").Value = "  -2.89%  "
$ws.Range("D15").Value = "57.766.96"
$ws.Range("E15").Value = "  -1.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.68"
$ws.Range("E16").Value = "  -2.13%  "
$ws.Range("E17").Value = "  -1.65%  "
$ws.Range("D18").Value = "2.449.36"
$ws.Range("E18").Value = "  -2.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.29"
$ws.Range("E19").Value = "  -3.61%  "
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "317.23"
$ws.Range("E21").Value = "  -1.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.10"
$ws.Range("E22").Value = "  -1.23%  "
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.64"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("E25").Value = "  -1.24%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -1.86%  "
$ws.Range("E28").Value = "  -2.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "174.46"
$ws.Range("E29").Value = "  +3.93%  "
$ws.Range("D30").Value = "0.0₃0735"
$ws.Range("E30").Value = "  -2.91%  "
$ws.Range("E31").Value = "  -1.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.12"
$ws.Range("E32").Value = "  -3.91%  "
$ws.Range("E33").Value = "  -5.59%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.80"
$ws.Range("E36").Value = "  -1.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.18"
$ws.Range("E37").Value = "  -6.34%  "
$ws.Range("E38").Value = "  -4.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.22"
$ws.Range("E39").Value = "  -0.42%  "
$ws.Range("E40").Value = "  -2.39%  "
$ws.Range("E41").Value = "  +2.87%  "
$ws.Range("E42").Value = "  -2.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "263.62"
$ws.Range("E43").Value = "  -5.53%  "
$ws.Range("E44").Value = "  -2.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.80"
$ws.Range("E45").Value = "  -4.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0923"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.85"
$ws.Range("E47").Value = "  -6.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0495"
$ws.Range("E48").Value = "  -1.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0211"
$ws.Range("E49").Value = "  -1.48%  "
$ws.Range("E50").Value = "  -4.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.31"
$ws.Range("E51").Value = "  -3.73%  "
